$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 144; this shifts the existing rows
# 144-190 down to 145-191 (preserving all of their data/formatting),
# matching the "add new weekly record" edit described by the diff.
$ws.Rows.Item(144).Insert()

# Populate the newly inserted row 144 with the new weekly record.
$ws.Range("A144").Value = 4
$ws.Range("B144").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C144").Value = "Los Lagos"
$ws.Range("D144").Value = 44559
$ws.Range("E144").Value = 10
$ws.Range("F144").Value = 100112021
$ws.Range("G144").Value = "Ají"
$ws.Range("H144").Value = "Inferno"
$ws.Range("I144").Value = "Primera"
$ws.Range("J144").Value = 140
$ws.Range("K144").Value = 24500
$ws.Range("L144").Value = 25000
$ws.Range("M144").Value = 24750
$ws.Range("N144").Value = "$/caja 12 kilos"
$ws.Range("O144").Value = "Región de Arica y Parinacota"
$ws.Range("P144").Value = 2062
$ws.Range("Q144").Value = 12
$ws.Range("R144").Value = "Hortaliza"
